$d = $word.ActiveDocument

# --- Locate the run that still holds the placeholder date text -------------
$searchText = "no Período de 00/00/000 a 00/00/000"

$targetPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.IndexOf($searchText) -ge 0) {
        $targetPara = $para
        break
    }
}

$pStart = $targetPara.Range.Start
$idxInPara = $targetPara.Range.Text.IndexOf($searchText)
$rngStart = $pStart + $idxInPara
$rngEnd = $rngStart + $searchText.Length

# Range used purely as a formatting donor (1 char is enough; it still points
# at "n" of the original "no Período..." run, so it carries the exact
# rFonts/b/bCs the new runs need to match).
$fmtRefStart = $rngStart

# --- Helper: insert `text` at `insertPos` as its own run, cloning the ------
# --- formatting of the donor range without disturbing neighbouring runs. ---
function Insert-FormattedPiece($insertPos, $text, $fmtRefStart) {
    $refRng = $d.Range($fmtRefStart, $fmtRefStart + 1)

    # Scratch paragraph at the very end of the document - isolated, so the
    # formatted-text round trip below can never coalesce with a neighbour.
    $docEnd = $d.Content.End - 1
    $scratchPara = $d.Paragraphs.Add($d.Range($docEnd, $docEnd))
    $scratchStart = $scratchPara.Range.Start

    $scratchTarget = $d.Range($scratchStart, $scratchStart)
    $scratchTarget.FormattedText = $refRng.FormattedText

    $scratchRunRng = $d.Range($scratchStart, $scratchStart + 1)
    $scratchRunRng.Text = $text

    $finalSrc = $d.Range($scratchStart, $scratchStart + $text.Length)
    $tgt = $d.Range($insertPos, $insertPos)
    $tgt.FormattedText = $finalSrc.FormattedText

    $scratchParaRange = $scratchPara.Range
    $delRng = $d.Range($scratchParaRange.Start, $scratchParaRange.End)
    $delRng.Delete()

    return $insertPos + $text.Length
}

# --- Insert the eight replacement pieces ahead of the old placeholder text -
$pieces = @("no Período de ", "01", "/", "08", "/", "2025", " a ", "31/08/2025")

$insertPos = $rngStart
foreach ($piece in $pieces) {
    $insertPos = Insert-FormattedPiece $insertPos $piece $fmtRefStart
}

# --- Remove the now-shifted-right original placeholder text ----------------
$totalInserted = 0
foreach ($piece in $pieces) { $totalInserted += $piece.Length }

$oldTextStart = $rngStart + $totalInserted
$oldTextEnd = $oldTextStart + $searchText.Length
$oldRng = $d.Range($oldTextStart, $oldTextEnd)
$oldRng.Delete()

Write-Host "Final paragraph text:" $targetPara.Range.Text
